$wb = $excel.ActiveWorkbook

# Add a new worksheet for the ValidateSidebar test case and place it
# right after Sheet1 (Excel's default Add() puts a new sheet before the
# active sheet, so move it explicitly into the 2nd position).
$new = $wb.Worksheets.Add()
$new.Name = "Sheet2"
$wb.Sheets("Sheet2").Move($null, $wb.Sheets("Sheet1"))

# Sheet handles track position, not identity, so re-fetch by name after
# the move above before writing to it.
$ws2 = $wb.Sheets("Sheet2")

# Populate the sidebar menu items.
$ws2.Range("A1").Value = "All Items"
$ws2.Range("A2").Value = "About"
$ws2.Range("A3").Value = "Logout"
$ws2.Range("A4").Value = "Reset App State"

# Select A4 and make Sheet2 the active sheet/tab.
$ws2.Activate()
$ws2.Range("A4").Select()
